# Applies the Review_466.docx edit: new title/date, restructured
# "Spurious Rewards" summary with Heading3 section headers, and a new
# arXiv link, replacing the old Transformer^2 review content.
$d = $word.ActiveDocument

# --- Paragraphs 1-7 already exist in the document; overwrite their
#     text/style in place to match the first 7 paragraphs of the target. ---
$p1 = $d.Paragraphs(1)
$p1.Range.Text = "המאמר היומי של יניב ומייק: 09.06.25" + [char]11 + "Spurious Rewards: Rethinking Training Signals in RLVR – Fast Overview"
$p1.Style = "Normal"

$p2 = $d.Paragraphs(2)
$p2.Range.Text = "המסר המרכזי במשפט אחד" + [char]11 + " גם תגמולים אקראיים או שגויים יכולים להביא לשיפור דרמטי ביכולות פתרון בעיות מתמטיות – אבל רק אם המודל כבר `"מכיר`" את הדרך מהפרה-טריינינג."
$p2.Style = "Normal"

$p3 = $d.Paragraphs(3)
$p3.Range.Text = "למה זה חשוב"
$p3.Style = "Heading3"

$p4 = $d.Paragraphs(4)
$p4.Range.Text = "למידה באמצעות חיזוקים עם תגמול ניתן לאימות (RL with Verifiable Rewards - RLVR) הפכה לשיטה מובילה לשפר יכולות חשיבה של מודלים גדולים. המאמר שואל שאלה פרובוקטיבית: האם אנחנו באמת צריכים תגמול מדויק? התשובה: לא תמיד."
$p4.Style = "Normal"

$p5 = $d.Paragraphs(5)
$p5.Range.Text = "מה עשו החוקרים"
$p5.Style = "Heading3"

$p6 = $d.Paragraphs(6)
$p6.Range.Text = "הם לקחו את המודל Qwen-2.5-Math ואימנו אותו על סט שאלות מתמטיקה עם חמש גרסאות שונות של תגמולים:"
$p6.Style = "Normal"

$p7 = $d.Paragraphs(7)
$p7.Range.Text = "תגמול אמיתי:  מודל מקבל נקודה רק אם התשובה נכונה."
$p7.Style = "Normal"

# --- Append the remaining new paragraphs (8-20) after paragraph 7,
#     each one created via InsertParagraphAfter chained off the previous. ---
$anchor = $d.Paragraphs(7)
$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs(8)
$newp.Range.Text = "תגמול לפי הצבעת רוב: המודל מייצר 64 תשובות, ומתגמל את התשובה השכיחה." + [char]11 + "תגמול פורמטי: אם התשובה כוללת ביטוי מתמטי (למשל \boxed{}), היא מתוגמלת, בלי קשר לנכונות."
$newp.Style = "Normal"
$anchor = $newp

$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs(9)
$newp.Range.Text = "תגמול אקראי: הטלת מטבע קובעת אם לתגמל."
$newp.Style = "Normal"
$anchor = $newp

$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs(10)
$newp.Range.Text = "תגמול הפוך: רק תשובות שגויות מקבלות נקודה."
$newp.Style = "Normal"
$anchor = $newp

$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs(11)
$newp.Range.Text = "במפתיע, כל אחד מהתגמולים הללו הצליח כמעט כמו תגמול אמיתי כלומר המודל השתפר דרמטית גם כש האות החיזוקי לא היה קשור כלל לתוצאה הנכונה."
$newp.Style = "Normal"
$anchor = $newp

$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs(12)
$newp.Range.Text = "ממצאים עיקריים"
$newp.Style = "Heading3"
$anchor = $newp

$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs(13)
$newp.Range.Text = "Qwen משתפר בכל תנאי: גם בלי תגמול נכון, המודל לומד לפתור בעיות טוב יותר. לעומת זאת, מודלים אחרים (כמו Llama3 ו־OLMo2) זקוקים לתגמול מדויק כדי להשתפר."
$newp.Style = "Normal"
$anchor = $newp

$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs(14)
$newp.Range.Text = "הגורם הסמוי: פתרון דרך קוד. Qwen כבר יודע לנסח פתרונות בפייתון מתוך הטקסט. אימון RLVR רק גורם לו לבחור באסטרטגיה הזו יותר ומביא לדיוק גבוה יותר."
$newp.Style = "Normal"
$anchor = $newp

$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs(15)
$newp.Range.Text = "שיפור בדיוק נובע ממעבר מ`"לשוני`" ל`"קוד`": בשאלות שבהן המודל התחיל לכתוב קוד בעקבות האימון, הדיוק קפץ בכמעט 26%."
$newp.Style = "Normal"
$anchor = $newp

$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs(16)
$newp.Range.Text = "אז למה תגמול אקראי עובד? האלגוריתם GRPO כולל קליפינג שמעדיף פעולות בסבירות גבוהה – כך שגם כשאין קשר לתוצאה, המודל לומד לחזק את ההתנהגות הדומיננטית שלו."
$newp.Style = "Normal"
$anchor = $newp

$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs(17)
$newp.Range.Text = "לא כל מודל נולד שווה: כשאין במודל נטייה מוקדמת לקוד, כמו ב־OLMo, אותו תגמול אקראי פשוט לא עובד." + [char]11
$newp.Style = "Normal"
$anchor = $newp

$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs(18)
$newp.Range.Text = "סיכום:"
$newp.Style = "Normal"
$anchor = $newp

$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs(19)
$newp.Range.Text = "המאמר מראה שלעיתים קרובות אימון RL לא מלמד כישורים חדשים, אלא מחלץ כישורים חבויים שהמודל כבר פיתח בפרה-טריינינג. לא תמיד צריך תגמול מדויק – אם המודל כבר `"מכיר`" את הדרך, מספיק לאותת לו לחזור אליה. עם זאת, זה לא נכון לכל מודל – יש כאלה שדורשים הנחיה מדויקת כדי להשתפר."
$newp.Style = "Normal"
$anchor = $newp

$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs(20)
$newp.Range.Text = "https://arxiv.org/abs/2412.07169 "
$newp.Style = "Normal"
$anchor = $newp

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
